$d = $word.ActiveDocument

# Find the list-paragraph that currently reads
# "Have level restart. Check with video or play game to see how level resets."
# and replace its runs with the new note about resetting the ghost/Pac-Man,
# matching how Word splits runs around proofing-flagged words
# ("StopAllCouroutine" and "pac") with w:proofErr spellStart/spellEnd markers.
$findRng = $d.Content
$findRng.Find.Execute("Have level restart.", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null

$p = $findRng.Paragraphs(1)
$rng = $p.Range

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00564ACD" w:rsidRDefault="00564ACD" w:rsidP="00DF5ECA"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>On reset, have ghost start again.</w:t></w:r><w:r><w:t xml:space="preserve"> Make sure to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StopAllCouroutine</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pac</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> man and ghost.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
